# TestFlow.xlsx update: add CBP_PXS_001 keyword-driven test row (px search / upax)
# and rename the two existing CBP test cases to use the shared OPEN_BROWSER /
# CBP_LOGIN / CLOSE_BROWSER keywords.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: extend Keyword6..Keyword15 across new columns H:P.
# Copy the formatting from the existing header cell (G1) before writing
# the new header labels so the new cells pick up the bold/centered style.
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 8))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 9))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 10))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 11))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 12))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 13))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 14))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 15))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 16))

$ws.Cells.Item(1, 8).Value  = "Keyword7"
$ws.Cells.Item(1, 9).Value  = "Keyword8"
$ws.Cells.Item(1, 10).Value = "Keyword9"
$ws.Cells.Item(1, 11).Value = "Keyword10"
$ws.Cells.Item(1, 12).Value = "Keyword11"
$ws.Cells.Item(1, 13).Value = "Keyword12"
$ws.Cells.Item(1, 14).Value = "Keyword13"
$ws.Cells.Item(1, 15).Value = "Keyword14"
$ws.Cells.Item(1, 16).Value = "Keyword15"

# --- Row 2: CBP001 / CBP_LOGIN / CLOSE_BROWSER ---
$ws.Cells.Item(2, 1).Value = "CBP001"
$ws.Cells.Item(2, 3).Value = "CBP_LOGIN"
$ws.Cells.Item(2, 4).Value = "CLOSE_BROWSER"
$ws.Cells.Item(2, 5).ClearContents()

# --- Row 3: CBP002 / CBP_LOGIN / CREATE_AND_FILL_1DAY_LOOKOUT / CLOSE_BROWSER ---
$ws.Cells.Item(3, 1).Value = "CBP002"
$ws.Cells.Item(3, 3).Value = "CBP_LOGIN"
$ws.Cells.Item(3, 4).Value = "CREATE_AND_FILL_1DAY_LOOKOUT"
$ws.Cells.Item(3, 5).Value = "CLOSE_BROWSER"
$ws.Cells.Item(3, 6).ClearContents()

# --- Row 4: CBP_PXS_001, the new px-search / upax event flow.
# Copy formatting from an existing data cell (G4) into the new columns
# H:P before filling in the keyword values.
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 8))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 9))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 10))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 11))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 12))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 13))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 14))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 15))
$ws.Cells.Item(4, 7).Copy($ws.Cells.Item(4, 16))

$ws.Cells.Item(4, 1).Value  = "CBP_PXS_001"
$ws.Cells.Item(4, 3).Value  = "CBP_LOGIN"
$ws.Cells.Item(4, 4).Value  = "NAVIGATE_TO_PERSON_SEARCH"
$ws.Cells.Item(4, 5).Value  = "SEARCH_PERSON"
$ws.Cells.Item(4, 6).Value  = "SELECT_PXSEARCH"
$ws.Cells.Item(4, 7).Value  = "CREATE_UPAX_EVENT_EXISTING"
$ws.Cells.Item(4, 8).Value  = "SWITCH_TO_NEW_TAB"
$ws.Cells.Item(4, 9).Value  = "SELECT_EXISTING_EVENT_TAB"
$ws.Cells.Item(4, 10).Value = "ENTER_EVENT_NUMBER_AND_SELECT"
$ws.Cells.Item(4, 11).Value = "SELECT_NEW_PERSON_TAB"
$ws.Cells.Item(4, 12).Value = "SET_PERSON_TYPE_AND_IMPORT"
$ws.Cells.Item(4, 13).Value = "SELECT_ASSOCIATED_PERSON"
$ws.Cells.Item(4, 14).Value = "VERIFY_SUBJECT"
$ws.Cells.Item(4, 15).Value = "DELETE_TRAVELER"
$ws.Cells.Item(4, 16).Value = "CLOSE_BROWSER"

# --- Row heights: row 2 now wraps to a single line (30), rows 3 & 4 grow
# to fit the extra wrapped keyword columns (75) ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 75

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection left below the data, matching the author's last interaction ---
$null = $ws.Range("A5:XFD6").Select()
